$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 233.83333
$ws.Range("I38").Value = 233.83333
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 701.49999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -329.49999
$ws.Range("N38").ClearContents()

$ws.Range("H64").Value = 76928270
$ws.Range("I64").Value = 5820
$ws.Range("J64").Value = 142861800
$ws.Range("K64").Value = 5820
$ws.Range("L64").Value = 142861800
$ws.Range("M64").Value = -5572
$ws.Range("N64").Value = -142862296

$ws.Range("H67").Value = 76928270
$ws.Range("I67").Value = 5820
$ws.Range("J67").Value = 142861800
$ws.Range("K67").Value = 5820
$ws.Range("L67").Value = 142861800
$ws.Range("M67").Value = -4962
$ws.Range("N67").Value = -142863516

$ws.Range("H94").Value = 1899.8
$ws.Range("I94").Value = 1899.8
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1899.8
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1448.8

$ws.Range("H98").Value = 3377.8572
$ws.Range("I98").Value = 3391.8518
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 3391.8518
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = -1893.8518
$ws.Range("N98").Value = -5996

$ws.Range("H111").Value = 1145.75
$ws.Range("I111").Value = 964.1429000000001
$ws.Range("J111").Value = 1400
$ws.Range("K111").Value = 2892.4287
$ws.Range("L111").Value = 4200
$ws.Range("M111").Value = 174.5712999999996
$ws.Range("N111").Value = -10334

$ws.Range("H116").Value = 10567.6
$ws.Range("I116").Value = 5096.7144
$ws.Range("J116").Value = 23333
$ws.Range("K116").Value = 5096.7144
$ws.Range("L116").Value = 23333
$ws.Range("M116").Value = -1654.7144
$ws.Range("N116").Value = -30217

$ws.Range("H122").Value = 3377.8572
$ws.Range("I122").Value = 3391.8518
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 10175.5554
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -7725.555399999999
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2004.5555
$ws.Range("I2").Value = 1596.6666
$ws.Range("J2").Value = 2208.5
$ws.Range("K2").Value = 1596.6666
$ws.Range("L2").Value = 2208.5
$ws.Range("M2").Value = -1483.6666
$ws.Range("N2").Value = -2434.5

$ws.Range("H5").Value = 191.25
$ws.Range("I5").Value = 191.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 191.25
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -79.25

$ws.Range("H32").Value = 1799.7826
$ws.Range("I32").Value = 1843.091
$ws.Range("J32").Value = 847
$ws.Range("K32").Value = 1843.091
$ws.Range("L32").Value = 847
$ws.Range("M32").Value = -1556.091
$ws.Range("N32").Value = -1421

$ws.Range("H61").Value = 3022.3845
$ws.Range("I61").Value = 2275.5386
$ws.Range("J61").Value = 3769.2307
$ws.Range("K61").Value = 2275.5386
$ws.Range("L61").Value = 3769.2307
$ws.Range("M61").Value = -2063.5386
$ws.Range("N61").Value = -4193.2307

$ws.Range("H74").Value = 66775.625
$ws.Range("I74").Value = 74249.12
$ws.Range("J74").Value = 3666.111
$ws.Range("K74").Value = 74249.12
$ws.Range("L74").Value = 3666.111
$ws.Range("M74").Value = -73375.12
$ws.Range("N74").Value = -5414.111

$ws.Range("H76").Value = 80000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 80000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 80000
$ws.Range("N76").Value = -80676

$ws.Range("H77").Value = 66775.625
$ws.Range("I77").Value = 74249.12
$ws.Range("J77").Value = 3666.111
$ws.Range("K77").Value = 371245.6
$ws.Range("L77").Value = 18330.555
$ws.Range("M77").Value = -366877.6
$ws.Range("N77").Value = -27066.555

$ws.Range("H79").Value = 80000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 80000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 80000
$ws.Range("N79").Value = -82340

$ws.Range("H116").Value = 2004.5555
$ws.Range("I116").Value = 1596.6666
$ws.Range("J116").Value = 2208.5
$ws.Range("K116").Value = 1596.6666
$ws.Range("L116").Value = 2208.5
$ws.Range("M116").Value = 697.3334
$ws.Range("N116").Value = -6796.5

$ws.Range("H136").Value = 3022.3845
$ws.Range("I136").Value = 2275.5386
$ws.Range("J136").Value = 3769.2307
$ws.Range("K136").Value = 6826.6158
$ws.Range("L136").Value = 11307.6921
$ws.Range("M136").Value = -4276.6158
$ws.Range("N136").Value = -16407.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2004.5555
$ws.Range("I3").Value = 1596.6666
$ws.Range("J3").Value = 2208.5
$ws.Range("K3").Value = 1596.6666
$ws.Range("L3").Value = 2208.5
$ws.Range("M3").Value = -1482.6666
$ws.Range("N3").Value = -2436.5

$ws.Range("H4").Value = 191.25
$ws.Range("I4").Value = 191.25
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 191.25
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -76.25

$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H55").Value = 23921.25
$ws.Range("I55").Value = 10505
$ws.Range("J55").Value = 28393.334
$ws.Range("K55").Value = 10505
$ws.Range("L55").Value = 28393.334
$ws.Range("M55").Value = -10232
$ws.Range("N55").Value = -28939.334

$ws.Range("H122").Value = 160000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 160000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 160000
$ws.Range("N122").Value = -169800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 297
$ws.Range("I7").Value = 328
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 328
$ws.Range("L7").Value = 80
$ws.Range("M7").Value = -215
$ws.Range("N7").Value = -306

$ws.Range("H22").Value = 2875
$ws.Range("I22").Value = 3500
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 3500
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -3150
$ws.Range("N22").Value = -1700

$ws.Range("H31").Value = 6950321.5
$ws.Range("I31").Value = 5079.7
$ws.Range("J31").Value = 15631874
$ws.Range("K31").Value = 5079.7
$ws.Range("L31").Value = 15631874
$ws.Range("M31").Value = -4784.7
$ws.Range("N31").Value = -15632464

$ws.Range("H34").Value = 6950321.5
$ws.Range("I34").Value = 5079.7
$ws.Range("J34").Value = 15631874
$ws.Range("K34").Value = 5079.7
$ws.Range("L34").Value = 15631874
$ws.Range("M34").Value = -4877.7
$ws.Range("N34").Value = -15632278

$ws.Range("H86").Value = 9331
$ws.Range("I86").Value = 7993
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 7993
$ws.Range("L86").Value = 10000
$ws.Range("M86").Value = -6870
$ws.Range("N86").Value = -12246

$ws.Range("H89").Value = 9331
$ws.Range("I89").Value = 7993
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 39965
$ws.Range("L89").Value = 50000
$ws.Range("M89").Value = -34349
$ws.Range("N89").Value = -61232

$ws.Range("H122").Value = 2677.96
$ws.Range("I122").Value = 1835.4
$ws.Range("J122").Value = 3941.8
$ws.Range("K122").Value = 5506.200000000001
$ws.Range("L122").Value = 11825.4
$ws.Range("M122").Value = -3056.200000000001
$ws.Range("N122").Value = -16725.4

$ws.Range("H132").Value = 37041030
$ws.Range("I132").Value = 5336.3335
$ws.Range("J132").Value = 55558880
$ws.Range("K132").Value = 16009.0005
$ws.Range("L132").Value = 166676640
$ws.Range("M132").Value = -13479.0005
$ws.Range("N132").Value = -166681700

$ws.Range("H134").Value = 3654.4583
$ws.Range("I134").Value = 3078.1667
$ws.Range("J134").Value = 5383.3335
$ws.Range("K134").Value = 9234.500100000001
$ws.Range("L134").Value = 16150.0005
$ws.Range("M134").Value = -6699.500100000001
$ws.Range("N134").Value = -21220.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1093.1395
$ws.Range("I2").Value = 497.5
$ws.Range("J2").Value = 1446.1111
$ws.Range("K2").Value = 2985
$ws.Range("L2").Value = 8676.6666
$ws.Range("M2").Value = -2872
$ws.Range("N2").Value = -8902.6666

$ws.Range("H6").Value = 98.5
$ws.Range("I6").Value = 98.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 295.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -182.5

$ws.Range("H60").Value = 742844.5600000001
$ws.Range("I60").Value = 1177365.9
$ws.Range("J60").Value = 4158.4
$ws.Range("K60").Value = 3532097.7
$ws.Range("L60").Value = 12475.2
$ws.Range("M60").Value = -3531846.7
$ws.Range("N60").Value = -12977.2

$ws.Range("H80").Value = 5350.5
$ws.Range("I80").Value = 702
$ws.Range("J80").Value = 9999
$ws.Range("K80").Value = 2106
$ws.Range("L80").Value = 29997
$ws.Range("M80").Value = -1170
$ws.Range("N80").Value = -31869

$ws.Range("H83").Value = 5350.5
$ws.Range("I83").Value = 702
$ws.Range("J83").Value = 9999
$ws.Range("K83").Value = 6318
$ws.Range("L83").Value = 89991
$ws.Range("M83").Value = -1638
$ws.Range("N83").Value = -99351

$ws.Range("H131").Value = 8144.7
$ws.Range("I131").Value = 8815.111000000001
$ws.Range("J131").Value = 2111
$ws.Range("K131").Value = 26445.333
$ws.Range("L131").Value = 6333
$ws.Range("M131").Value = -21405.333
$ws.Range("N131").Value = -16413

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 169.85715
$ws.Range("I2").Value = 72.75
$ws.Range("J2").Value = 299.33334
$ws.Range("K2").Value = 72.75
$ws.Range("L2").Value = 299.33334
$ws.Range("M2").Value = 40.25
$ws.Range("N2").Value = -525.33334

$ws.Range("H52").Value = 50001
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 50001
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 50001
$ws.Range("N52").Value = -50519

$ws.Range("H107").Value = 4684.7896
$ws.Range("I107").Value = 419.5
$ws.Range("J107").Value = 11996.714
$ws.Range("K107").Value = 419.5
$ws.Range("L107").Value = 11996.714
$ws.Range("M107").Value = 1500.5
$ws.Range("N107").Value = -15836.714

$ws.Range("H126").Value = 13012.2
$ws.Range("I126").Value = 6537
$ws.Range("J126").Value = 22725
$ws.Range("K126").Value = 19611
$ws.Range("L126").Value = 68175
$ws.Range("M126").Value = -17141
$ws.Range("N126").Value = -73115

$ws.Range("H136").Value = 31355
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 31355
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 94065
$ws.Range("N136").Value = -99165

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1457.3572
$ws.Range("I16").Value = 1461.7693
$ws.Range("J16").Value = 1400
$ws.Range("K16").Value = 1461.7693
$ws.Range("L16").Value = 1400
$ws.Range("M16").Value = -1291.7693
$ws.Range("N16").Value = -1740

$ws.Range("H22").Value = 1949.2
$ws.Range("I22").Value = 2249
$ws.Range("J22").Value = 1499.5
$ws.Range("K22").Value = 2249
$ws.Range("L22").Value = 1499.5
$ws.Range("M22").Value = -1954
$ws.Range("N22").Value = -2089.5

$ws.Range("H27").Value = 1949.2
$ws.Range("I27").Value = 2249
$ws.Range("J27").Value = 1499.5
$ws.Range("K27").Value = 2249
$ws.Range("L27").Value = 1499.5
$ws.Range("M27").Value = -2142
$ws.Range("N27").Value = -1713.5

$ws.Range("H46").Value = 3716.6667
$ws.Range("I46").Value = 2825
$ws.Range("J46").Value = 5500
$ws.Range("K46").Value = 2825
$ws.Range("L46").Value = 5500
$ws.Range("M46").Value = -2637
$ws.Range("N46").Value = -5876

$ws.Range("H132").Value = 2829.2292
$ws.Range("I132").Value = 2384.6775
$ws.Range("J132").Value = 3639.8823
$ws.Range("K132").Value = 7154.032499999999
$ws.Range("L132").Value = 10919.6469
$ws.Range("M132").Value = -4624.032499999999
$ws.Range("N132").Value = -15979.6469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 34608.8
$ws.Range("I54").Value = 14535
$ws.Range("J54").Value = 47991.332
$ws.Range("K54").Value = 14535
$ws.Range("L54").Value = 47991.332
$ws.Range("M54").Value = -14015
$ws.Range("N54").Value = -49031.332

$ws.Range("H113").Value = 660.9677
$ws.Range("I113").Value = 590.25
$ws.Range("J113").Value = 789.5454999999999
$ws.Range("K113").Value = 1770.75
$ws.Range("L113").Value = 2368.6365
$ws.Range("M113").Value = 399.25
$ws.Range("N113").Value = -6708.6365

$ws.Range("H136").Value = 13507.875
$ws.Range("I136").Value = 8965.912
$ws.Range("J136").Value = 90721.25
$ws.Range("K136").Value = 26897.736
$ws.Range("L136").Value = 272163.75
$ws.Range("M136").Value = -24347.736
$ws.Range("N136").Value = -277263.75
